$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: I1 = "I0", J1 = "IF"
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from H1 (bold, bordered, centered header style) to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows 2-17 for columns I (I0) and J (IF)
$data = @{
    2  = @(1, 6)
    3  = @(1, 5)
    4  = @(1, 6)
    5  = @(1, 6)
    6  = @(1, 4)
    7  = @(1, 6)
    8  = @(1, 4)
    9  = @(5, 6)
    10 = @(4, 5)
    11 = @(7, 8)
    12 = @(8, 8)
    13 = @(8, 8)
    14 = @(4, 6)
    15 = @(6, 7)
    16 = @(8, 8)
    17 = @(9, 9)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
